$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: keep integer-style number format (style reused from existing table)
$ws.Range("A2").NumberFormat = "0"

# C2: replace the old text CNPJ (shared string) with the new numeric CNPJ value
$ws.Range("C2").Value = 11756894000160
$ws.Range("C2").NumberFormat = "0"
$ws.Range("C2").Interior.Color = 0
$ws.Range("C2").Font.Color = 16119285
$ws.Range("C2").Font.Size = 12
$ws.Range("C2").Font.Name = "Arial"

# Row height for row 2 (15.75pt)
$ws.Rows.Item(2).RowHeight = 15.75

# Column C width matches column A/B now
$ws.Columns.Item(3).ColumnWidth = 22.85546875

# Selection moves to C2
$ws.Range("C2").Select()
